$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the KS Tests functionality was verified (PASSED)
$passedRows = @(2, 3, 4, 5, 8, 9, 11, 12, 14, 15)

# Rows where the test failed (FAILED)
$failedRows = @(6, 10, 13, 16)

$passedActual = "Ks Tests functionality verified"
$failedActual = "Test failed - actual behavior did not match expected result"

# Row 7 already carries the exact "PASSED" formatting (green fill) that the
# other passing rows need to reuse, so copy its Actual Result / Test Status
# cell formatting onto each passing row instead of re-deriving a new style.
$passedFormatSource = $ws.Range("I7")

foreach ($r in $passedRows) {
    $ws.Range("H$r").Value = $passedActual
    $ws.Range("I$r").Value = "PASSED"

    $passedFormatSource.Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
}

# Establish the new "FAILED" style (red fill) on the first failing row, then
# reuse that exact formatting for the remaining failing rows so only one new
# fill/style gets introduced.
$firstFailedRow = $failedRows[0]
$ws.Range("H$firstFailedRow").Value = $failedActual
$ws.Range("I$firstFailedRow").Value = "FAILED"
$ws.Range("I$firstFailedRow").Interior.Color = 13551615

$failedFormatSource = $ws.Range("I$firstFailedRow")

foreach ($r in $failedRows) {
    if ($r -eq $firstFailedRow) { continue }
    $ws.Range("H$r").Value = $failedActual
    $ws.Range("I$r").Value = "FAILED"

    $failedFormatSource.Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
}
